$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# The merge-field placeholder "caseManagementLocation.external_short_name"
# (originally split over several runs: ".", "external", "_s", "hort",
# "_n", "ame") is replaced by "caseManagementLocation.venue_name".
# Only the ".external_short_name" portion (everything after the
# "caseManagementLocation" run) needs to change.
$d.Content.Find.Execute(".external_short_name", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ".venue_name", 2)

# --- Change 2 -------------------------------------------------------
# The heading "Witness statements" gets re-typed by the author, which
# causes Word to split it into two runs ("Witness " / "statements")
# flagged by the grammar checker (<w:proofErr w:type="gramStart"/> ...
# <w:proofErr w:type="gramEnd"/>). "statements" is unique in the
# document, so find it and force a run split at its boundary by
# toggling Bold off/on (text itself is unchanged).
$rng = $d.Content
$rng.Find.Execute("statements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 0
$rng.Font.Bold = 1
